# Updates cryptocurrency price/volume figures (and restores the original
# ARBITRUM/FraxShare row order) per the Nov 29 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new text value. NumberFormat is forced to
# "@" (Text) before the write and the cell style is reset to Normal right
# after, so numeric-looking strings (e.g. "226.75", "1.00") are preserved
# verbatim as text instead of being coerced into Excel numbers.
$updates = [ordered]@{
    "D2" = '37.720.11'
    "E2" = '  -1.26%  '
    "D3" = '2.024.48'
    "E3" = '  -1.73%  '
    "E4" = '  +0.31%  '
    "D5" = '226.75'
    "E5" = '  -1.69%  '
    "D6" = '0.612'
    "E6" = '  -0.63%  '
    "D7" = '59.71'
    "E7" = '  +1.91%  '
    "E8" = '  +0.16%  '
    "D9" = '0.384'
    "E9" = '  -1.09%  '
    "D10" = '0.0810'
    "E10" = '  +0.18%  '
    "E11" = '  -0.08%  '
    "D12" = '14.53'
    "E12" = '  -1.01%  '
    "D13" = '2.328.68'
    "E13" = '  -1.45%  '
    "D14" = '20.89'
    "E14" = '  +0.87%  '
    "D15" = '0.755'
    "E15" = '  -0.15%  '
    "D16" = '5.17'
    "E16" = '  -2.39%  '
    "D17" = '2.025.95'
    "E17" = '  -1.55%  '
    "D18" = '37.670.56'
    "E18" = '  -1.06%  '
    "D19" = '6.02'
    "E19" = '  -1.99%  '
    "D20" = '69.62'
    "E20" = '  -0.43%  '
    "D21" = '0.0₃0820'
    "E21" = '  -1.44%  '
    "D22" = '224.41'
    "E23" = '  -0.01%  '
    "E24" = '  -2.54%  '
    "D25" = '2.19'
    "E25" = '  -2.80%  '
    "D26" = '9.24'
    "E26" = '  -1.06%  '
    "D27" = '165.15'
    "E27" = '  -0.70%  '
    "E28" = '  -4.02%  '
    "D29" = '18.88'
    "E29" = '  -1.08%  '
    "D30" = '1.28'
    "E30" = '  -5.90%  '
    "E31" = '  +0.64%  '
    "D32" = '4.42'
    "E32" = '  -2.99%  '
    "E33" = '  +4.35%  '
    "D34" = '4.48'
    "E34" = '  -3.14%  '
    "D35" = '0.0600'
    "E35" = '  -2.48%  '
    "D36" = '6.29'
    "E36" = '  +4.11%  '
    "D37" = '2.24'
    "E37" = '  -4.66%  '
    "D38" = '3.24'
    "E38" = '  -2.72%  '
    "E39" = '  +0.34%  '
    "D40" = '1.533.73'
    "E40" = '  +3.49%  '
    "D41" = '0.0216'
    "E41" = '  -1.37%  '
    "D42" = '96.52'
    "E42" = '  -1.99%  '
    "D43" = '16.76'
    "E43" = '  -0.94%  '
    "E44" = '  -0.47%  '
    "D45" = '0.0916'
    "E45" = '  -3.24%  '
    "D46" = '1.10'
    "E46" = '  -1.98%  '
    "E47" = '  -4.49%  '
    "B48" = 'FraxShare'
    "C48" = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    "D48" = '7.14'
    "E48" = '  +0.62%  '
    "B49" = 'ARBITRUM'
    "C49" = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    "D49" = '1.00'
    "E49" = '  -1.95%  '
    "E50" = '  -0.60%  '
    "D51" = '2.217.55'
    "E51" = '  -1.47%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
